$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the Test Case description text (B1)
$ws.Range("B1").Value = "Test Case: Testing to see if creating assessments reaches data base with no errors"

# Update existing Step 1 / Step 2 text (rows 2-3)
$ws.Range("C2").Value = "Step 1: While logged out go to the ""create assessment"" page"
$ws.Range("D2").Value = "I am returned to the login page"

$ws.Range("C3").Value = "Step 2: Login as a user with the appropriate role and go to the create assessment page"
$ws.Range("D3").Value = "I am redirected to the create assessment page"

# Add new Step 3-6 rows
$ws.Range("C4").Value = "Step 3: Fill out an assessment for an employee on any of your teams"
$ws.Range("D4").Value = "A new assessment will be added to the database for that employee"

$ws.Range("C5").Value = "Step 4: Fill out a assessment for an employee that is not in any of your teams"
$ws.Range("D5").Value = "I am denied access to this"

$ws.Range("C6").Value = "Step 5: Leave Certain fields that are required blank"
$ws.Range("D6").Value = "An error text shows saying that some of the fields have been left blank"

$ws.Range("C7").Value = "Step 6: Create an assessment for the current user logged in"
$ws.Range("D7").Value = "I am denied access to this"

# Apply the same formatting (wrap text, vertical top alignment) used by the other step cells
$ws.Range("C4:D7").WrapText = $true
$ws.Range("C4:D7").VerticalAlignment = -4160

# Update selection to match target workbook state
$ws.Range("A1:F7").Select()
